$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44972
$ws.Range("M2").Value = 140
$ws.Range("N2").Value = 27000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 27429
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1524
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44671
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 29000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 29500
$ws.Range("Q3").Value = "$/caja 20 kilos"
$ws.Range("S3").Value = 1475
$ws.Range("T3").Value = 20

# Row 4
$ws.Range("D4").Value = 45014
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("Q4").Value = "$/caja 20 kilos"
$ws.Range("S4").Value = 1225
$ws.Range("T4").Value = 20

# Row 6
$ws.Range("D6").Value = 45007
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 27000
$ws.Range("O6").Value = 28000
$ws.Range("P6").Value = 27500
$ws.Range("S6").Value = 1375

# Row 7
$ws.Range("D7").Value = 45028
$ws.Range("L7").Value = "Segunda"
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("S7").Value = 1075

# Row 8
$ws.Range("D8").Value = 44664
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 29000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 29500
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("S8").Value = 1639
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44679
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 29000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 29500
$ws.Range("S9").Value = 1475

# Row 10
$ws.Range("D10").Value = 44679
$ws.Range("L10").Value = "Tercera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 25000
$ws.Range("P10").Value = 24500
$ws.Range("S10").Value = 1225

# Row 11
$ws.Range("D11").Value = 44979

# Row 12
$ws.Range("D12").Value = 44636
$ws.Range("L12").Value = "Primera"
$ws.Range("N12").Value = 29000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 29500
$ws.Range("S12").Value = 1475

# Row 13
$ws.Range("D13").Value = 44965
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 34000
$ws.Range("O13").Value = 35000
$ws.Range("P13").Value = 34600
$ws.Range("S13").Value = 1922

# Row 14
$ws.Range("D14").Value = 44965
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 120
$ws.Range("N14").Value = 32000
$ws.Range("O14").Value = 33000
$ws.Range("P14").Value = 32333
$ws.Range("Q14").Value = "$/caja 18 kilos"
$ws.Range("S14").Value = 1796
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("D15").Value = 44993
$ws.Range("M15").Value = 130
$ws.Range("N15").Value = 25000
$ws.Range("O15").Value = 26000
$ws.Range("P15").Value = 25462
$ws.Range("S15").Value = 1273

# Row 16
$ws.Range("D16").Value = 45021
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 22500
$ws.Range("Q16").Value = "$/caja 20 kilos"
$ws.Range("R16").Value = "Región de Coquimbo"
$ws.Range("S16").Value = 1125
$ws.Range("T16").Value = 20

# Row 17
$ws.Range("D17").Value = 44650
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = 31000
$ws.Range("O17").Value = 32000
$ws.Range("P17").Value = 31500
$ws.Range("S17").Value = 1575

# Row 18
$ws.Range("D18").Value = 44650
